# Add a new batch of vocabulary entries (rows 128-144) to Sheet1, and
# update the saved cursor/selection position to match the author's
# final viewport (B145).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New word pairs (English / Chinese), with an occasional "day N" marker
# in column C, following the existing layout of the sheet.
$ws.Range("A128").Value = "privilege"
$ws.Range("B128").Value = "特权"
$ws.Range("C128").Value = "day 4"

$ws.Range("A129").Value = "invalid"
$ws.Range("B129").Value = "无效的"

$ws.Range("A130").Value = "module"
$ws.Range("B130").Value = "模块"

$ws.Range("A131").Value = "interface"
$ws.Range("B131").Value = "接口"

$ws.Range("A132").Value = "peripheral"
$ws.Range("B132").Value = "外围设备"

$ws.Range("A133").Value = "latency"
$ws.Range("B133").Value = "延迟"

$ws.Range("A134").Value = "protocal"
$ws.Range("B134").Value = "原型的"

$ws.Range("A135").Value = "integrated"
$ws.Range("B135").Value = "综合"

$ws.Range("A136").Value = "semiconductor"
$ws.Range("B136").Value = "半导体"

$ws.Range("A137").Value = "robust"
$ws.Range("B137").Value = "健壮的"

$ws.Range("A138").Value = "verify"
$ws.Range("B138").Value = "验证"

$ws.Range("A139").Value = "serial"
$ws.Range("B139").Value = "连载"

$ws.Range("A140").Value = "aerospace"
$ws.Range("B140").Value = "航空航天"

$ws.Range("A141").Value = "duplex"
$ws.Range("B141").Value = "复式"

$ws.Range("A142").Value = "simultaneously"
$ws.Range("B142").Value = "同时"

$ws.Range("A143").Value = "simplex"
$ws.Range("B143").Value = "单纯形"

$ws.Range("A144").Value = "simultaneously"
$ws.Range("B144").Value = "同时"

# Move the selection/cursor to match the author's final position.
$ws.Range("B145").Select()
